$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.969.15"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -2.66%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'3.364.58"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -2.20%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Value = "'  +0.04%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'567.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -1.74%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'147.82"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +0.44%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("E7").Value = "'  +0.04%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.481"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +0.47%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("E9").Value = "'  +0.44%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("E10").Value = "'  -0.95%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'0.416"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +2.35%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'3.949.89"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -1.96%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("E13").Value = "'  +0.55%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'27.97"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.95%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'3.363.75"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -2.40%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("E16").Value = "'  -1.18%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'61.023.49"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -2.69%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("E18").Value = "'  -1.19%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'14.49"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -0.77%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'8.95"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -1.28%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'374.97"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -2.94%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("E22").Value = "'  +0.17%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'75.43"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.63%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("E24").Value = "'  +0.05%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'3.506.18"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -2.31%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("E26").Value = "'  -6.52%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("E27").Value = "'  -3.46%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("B28").Value = 'RenderToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D28").Value = "'7.42"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -2.44%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("B29").Value = 'Binance-PegBSC-USD'
$ws.Range("C29").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D29").Value = "'0.990"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -1.09%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("E30").Value = "'  -0.75%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'1.00"
$ws.Range("D31").Style = "Normal"

$ws.Range("D32").Value = "'7.71"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -3.17%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'22.85"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -1.62%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("E34").Value = "'  -3.72%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("E35").Value = "'  +1.06%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = "'169.82"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -0.24%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("E37").Value = "'  -3.84%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = "'6.81"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -2.43%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'29.16"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -8.31%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'3.402.51"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -2.12%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'0.0756"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -2.41%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("E42").Value = "'  -3.47%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("E43").Value = "'  -0.92%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("E45").Value = "'  -5.06%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'2.492.09"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -2.64%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'22.70"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +0.28%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("E48").Value = "'  -3.01%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("E49").Value = "'  +0.07%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("E50").Value = "'  -2.07%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'0.816"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +0.30%  "
$ws.Range("E51").Style = "Normal"

